$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "Study" header (column A) to lowercase "study" ahead of
# reworking the forest plot to show standardized mean differences.
$ws.Range("A1").Value = "study"

# Reset the lingering stale selection (was I34) left over in the sheet view.
[void]$ws.Range("A1").Select()
